$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")

$ws.Range("A16").Value = "sandro"

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "1238192312"

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "91823918239"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1928391823"

$ws.Range("F16").Value = "hduhwdq@djuqwhduqh"

$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "1283128312"

$ws.Range("H16").Value = "rua das alamedas"

$ws.Range("A1").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("G16").PasteSpecial(-4122)
